# New PO forecast model
# Updates "Weekly Quantity", "Monthly Trend" and "PO Forecast" sheets with
# refreshed forecast data: appends newly observed weekly/monthly rows and
# replaces the PO Forecast series with the output of the new model.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: Weekly Quantity  -> append rows 51-54
# ---------------------------------------------------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")

$weeklyNewRows = @(
    @(51, 45662.99999999999, 79),
    @(52, 45669.99999999999, 2),
    @(53, 45683.99999999999, 22),
    @(54, 45690.99999999999, 1)
)

foreach ($row in $weeklyNewRows) {
    $r = $row[0]
    $cellA = $wsWeekly.Cells.Item($r, 1)
    $cellA.Value = $row[1]
    $cellA.NumberFormat = $wsWeekly.Cells.Item($r - 1, 1).NumberFormat
    $wsWeekly.Cells.Item($r, 2).Value = $row[2]
}

# ---------------------------------------------------------------------
# Sheet: Monthly Trend -> append row 17
# ---------------------------------------------------------------------
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

$monthlyNewRow = 17
$monthlyNewCellA = $wsMonthly.Cells.Item($monthlyNewRow, 1)
$monthlyNewCellA.Value = 45688.99999999999
$monthlyNewCellA.NumberFormat = $wsMonthly.Cells.Item($monthlyNewRow - 1, 1).NumberFormat
$wsMonthly.Cells.Item($monthlyNewRow, 2).Value = 104

# ---------------------------------------------------------------------
# Sheet: PO Forecast -> replace forecast quantities (B2:B50), replace
# dates+quantities for rows 51-58, and append new forecast rows 59-62
# ---------------------------------------------------------------------
$wsForecast = $wb.Worksheets.Item("PO Forecast")

# New forecast quantities for existing dates (rows 2-50, dates unchanged)
$forecastQtyUpdates = @(
    @(2, 256), @(3, 244), @(4, 226), @(5, 159), @(6, 111),
    @(7, 141), @(8, 213), @(9, 228), @(10, 157), @(11, 92),
    @(12, 361), @(13, 0), @(14, 91), @(15, 71), @(16, 5),
    @(17, 30), @(18, 99), @(19, 99), @(20, 37), @(21, 26),
    @(22, 114), @(23, 213), @(24, 215), @(25, 129), @(26, 67),
    @(27, 94), @(28, 178), @(29, 162), @(30, 90), @(31, 32),
    @(32, 47), @(33, 119), @(34, 151), @(35, 91), @(36, 15),
    @(37, 34), @(38, 143), @(39, 209), @(40, 142), @(41, 17),
    @(42, 0), @(43, 44), @(44, 82), @(45, 20), @(46, 6),
    @(47, 673), @(48, 886), @(49, 75), @(50, 168)
)

foreach ($u in $forecastQtyUpdates) {
    $wsForecast.Cells.Item($u[0], 2).Value = $u[1]
}

# Rows 51-58: both the date (A) and quantity (B) are replaced
$forecastReplacedRows = @(
    @(51, 45662.99999999999, 100),
    @(52, 45669.99999999999, 80),
    @(53, 45683.99999999999, 75),
    @(54, 45690.99999999999, 20),
    @(55, 45697.99999999999, 0),
    @(56, 45704.99999999999, 0),
    @(57, 45711.99999999999, 24),
    @(58, 45718.99999999999, 69)
)

foreach ($row in $forecastReplacedRows) {
    $r = $row[0]
    $wsForecast.Cells.Item($r, 1).Value = $row[1]
    $wsForecast.Cells.Item($r, 2).Value = $row[2]
}

# New forecast rows 59-62 (appended)
$forecastNewRows = @(
    @(59, 45725.99999999999, 24),
    @(60, 45732.99999999999, 0),
    @(61, 45739.99999999999, 0),
    @(62, 45746.99999999999, 40)
)

foreach ($row in $forecastNewRows) {
    $r = $row[0]
    $cellA = $wsForecast.Cells.Item($r, 1)
    $cellA.Value = $row[1]
    $cellA.NumberFormat = $wsForecast.Cells.Item($r - 1, 1).NumberFormat
    $wsForecast.Cells.Item($r, 2).Value = $row[2]
}
